$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for season record: Wins, Losses, Ties
# Copy the existing header formatting (bold, centered, bordered) from AC1
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record values (Wins/Losses/Ties) for every player row
$lastRow = 51
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 64
    $ws.Cells.Item($r, 31).Value = 98
    $ws.Cells.Item($r, 32).Value = 0
}
